$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell: C1 = "~" ---
$ws.Range("C1").Value = "~"

# --- Alignment style #1: horizontal=center, vertical=center (applied to C1 "~") ---
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4108

# --- Alignment style #2: horizontal=left, vertical=center (applied to the header cells) ---
# Build the style cleanly on A1 first (single distinct cell -> no orphan xf),
# then propagate it to the rest of the header cells via copy/paste-special (format only)
# so every cell lands on the very same style index instead of minting new ones.
$ws.Range("A1").HorizontalAlignment = -4131
$ws.Range("A1").VerticalAlignment = -4108

$ws.Range("A1").Copy()
$ws.Range("A2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column A width: 12.75 -> 16 characters ---
$ws.Columns.Item(1).ColumnWidth = 15.285714285714286

# --- Selection left at C9, matching the final saved view ---
$ws.Range("C9").Select()
